$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 496.57144
